$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new columns (F, G) are being added for "height" and "weight" data.
# Copy the header formatting (bold font + border + centered alignment) from
# the existing E1 header into the two new header cells.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# The "height"/"weight" headers take over columns E/F, and the old
# "fantasy points" header (and its data) shifts right into column G.
$ws.Range("E1").Value2 = "height"
$ws.Range("F1").Value2 = "weight"
$ws.Range("G1").Value2 = "fantasy points"

for ($r = 2; $r -le 11; $r++) {
    $oldFantasy = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value2 = 6.333333333333333
    $ws.Cells.Item($r, 6).Value2 = 233
    $ws.Cells.Item($r, 7).Value2 = $oldFantasy
}
